$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D2 previously held "TENERA SHEER SHADINGS CARERING"; it now holds the new
# "CONTINUUM" variant (adds a new shared string).
$ws.Range("D2").Value = "TENERA SHEER SHADINGS CONTINUUM"

# The old D2 text is preserved by moving it into the new column F2, and G2
# duplicates E2's value - both new cells adopt the same number/font style as
# the existing D2/E2 cells (cellXfs index 2: Consolas 7pt). Copying the
# format from E2 (instead of setting font properties ourselves) lets Excel
# reuse the existing style definition rather than minting new font/xf
# entries.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("F2").Value = "TENERA SHEER SHADINGS CARERING"
$ws.Range("G2").Value = "ANOTONIA 2 1/2"""

# Selection moves from D3 to the whole of row 3 (A3:XFD3).
$ws.Rows("3").Select() | Out-Null
